$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tool")

# --- Header row: drop the old person-focused columns, replace with the
#     BDF-tool-focused schema (developer_team, technical_area, id, name, description) ---
$ws.Range("A1").Value = "developer_team"
$ws.Range("B1").Value = "technical_area"
$ws.Range("C1").Value = "id"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "description"

# Old sheet had 7 columns (A:G); new schema only needs 5 (A:E), so clear the
# now-unused trailing columns. This also shrinks the sheet's used-range
# dimension back down to A1:E1.
$ws.Range("F1:G1").ClearContents()

# --- Data validations: drop the old vital_status (column D) list, add the
#     new developer_team (column A) and technical_area (column B) lists ---
$ws.Range("D2:D1048576").Validation.Delete()

$ws.Range("A2:A1048576").Validation.Add(3, 1, 1, '"ASKEM_NEU,ASKEM_MIT,ASKEM_NYU,ASKEM_Jataware,Netrias,N3C,BDC,CRA,DNAHIVE,HMS,Stanford,UAB,ICF,SageBio,Insilicom"')
$ws.Range("B2:B1048576").Validation.Add(3, 1, 1, '"TA1,TA2,TA3"')
